$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B23: audio text now mentions "next" audio (quoted)
$ws.Range("B23").Value = 'If true, plays "next" audio.  If false, does not.'

# Insert a fresh row at position 26 so the new B26 cell inherits the
# row-23..25 style (s="2") the same way Excel does when you insert a row
# above an existing formatted one. This pushes the old row 26 down to 27.
$ws.Rows("26:26").Insert(-4121)

# Re-populate row 26 with its original A-column text plus the new B-column value
$ws.Range("A26").Value = "runCountdownTimer() interval is done.  Repeat."
$ws.Range("B26").Value = "timerUI() displays 0 seconds."

# Remove the now-duplicated old row (shifted down to 27)
$ws.Rows("27:27").Delete()

# Update the selected cell in the sheet view to B27
$ws.Range("B27").Select()
